# Updates the cryptocurrency Price (column D) and Volume(1h) (column E)
# figures on Sheet1 to the latest scraped snapshot.
#
# Each target cell currently holds its number formatted as plain text
# (e.g. "330.05", "2.78%"), so we write the new text via Formula with a
# leading apostrophe (forces text entry, avoiding Excel auto-converting
# a numeric-looking / percent-looking string into a real number or an
# actual percentage value) and then reset the cell style back to
# "Normal" so no stray quote-prefix style lingers on the cell.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $newText) {
    $rng = $ws.Range($cellRef)
    $rng.Formula = "'" + $newText
    $rng.Style = "Normal"
}


Set-TextValue "D2" "329.65"
Set-TextValue "D3" "45.23"
Set-TextValue "E3" "2.34%"
Set-TextValue "D4" "5.491"
Set-TextValue "E4" "0.20%"
Set-TextValue "D5" "0.08424"
Set-TextValue "E5" "4.73%"
Set-TextValue "D6" "2.041"
Set-TextValue "E6" "0.77%"
Set-TextValue "D7" "0.9807"
Set-TextValue "E7" "2.99%"
Set-TextValue "D9" "0.1112"
Set-TextValue "E9" "-0.40%"
Set-TextValue "D10" "0.1916"
Set-TextValue "E10" "2.40%"
Set-TextValue "D11" "9.483"
Set-TextValue "E11" "-7.26%"
Set-TextValue "D12" "0.09729"
Set-TextValue "E12" "-1.64%"
Set-TextValue "D13" "0.04693"
Set-TextValue "E13" "-1.34%"
Set-TextValue "D14" "0.1060"
Set-TextValue "E14" "-0.31%"
Set-TextValue "E15" "2.45%"
Set-TextValue "D16" "0.04187"
Set-TextValue "E16" "2.32%"
Set-TextValue "D17" "0.006000"
Set-TextValue "E17" "4.36%"
Set-TextValue "D18" "3.386"
Set-TextValue "E18" "0.28%"
Set-TextValue "E19" "0.90%"
Set-TextValue "D20" "0.3354"
Set-TextValue "E20" "-1.60%"
Set-TextValue "D21" "0.1358"
Set-TextValue "E21" "-3.00%"
Set-TextValue "D22" "0.2553"
Set-TextValue "E22" "-1.06%"
Set-TextValue "E23" "-0.51%"
Set-TextValue "D24" "0.004446"
Set-TextValue "E24" "2.35%"
Set-TextValue "D25" "0.0001305"
Set-TextValue "E25" "4.09%"
Set-TextValue "E26" "-20.23%"
Set-TextValue "D38" "0.02716"
Set-TextValue "E38" "5.22%"
Set-TextValue "D39" "0.05655"
Set-TextValue "E39" "-0.01%"
Set-TextValue "D40" "0.007861"
Set-TextValue "E40" "1.58%"
Set-TextValue "D41" "0.1431"
Set-TextValue "E41" "2.38%"
Set-TextValue "D42" "0.007417"
Set-TextValue "E42" "0.78%"
Set-TextValue "D43" "0.002124"
Set-TextValue "E43" "5.52%"
Set-TextValue "D44" "0.008621"
Set-TextValue "E44" "1.03%"
Set-TextValue "D45" "0.3385"
Set-TextValue "D46" "0.00006895"
Set-TextValue "E46" "-2.96%"
Set-TextValue "D47" "0.00000000754"
Set-TextValue "E47" "0.27%"
Set-TextValue "E48" "0.18%"
Set-TextValue "D49" "0.003488"
Set-TextValue "E49" "-0.71%"
Set-TextValue "D50" "0.003549"
Set-TextValue "E50" "1.33%"
Set-TextValue "D51" "0.00002111"
Set-TextValue "E51" "0.27%"
